$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Country", "Practice Area", "Link", "Manager ", "Specialism", "Role", "Firm")
$col = 4
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$ws.Columns.ClearOutline()

$ws.Range("A2:J256").Select()
